$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.007.77'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '2.643.50'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '530.93'
$c.ClearFormats()
$ws.Range("E5").Value = '  +3.99%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '155.54'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  +0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.593'
$c.ClearFormats()
$ws.Range("E8").Value = '  +0.91%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  +4.81%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '3.106.97'
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").Value = '61.024.78'
$ws.Range("E14").Value = '  +1.17%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '21.96'
$c.ClearFormats()
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("D17").Value = '2.647.05'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("E18").Value = '  +0.22%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '354.41'
$c.ClearFormats()
$ws.Range("E19").Value = '  +1.17%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '10.65'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.63%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.22'
$c.ClearFormats()
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("E22").Value = '  +0.31%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '61.73'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("E25").Value = '  +1.50%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '0.0₃0860'
$ws.Range("E27").Value = '  +2.55%  '
$ws.Range("E28").Value = '  +0.52%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +4.14%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.15'
$c.ClearFormats()
$ws.Range("E31").Value = '  +7.02%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '19.51'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.58%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '150.33'
$c.ClearFormats()
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  +3.81%  '
$ws.Range("E35").Value = '  +1.79%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.922'
$c.ClearFormats()
$ws.Range("E36").Value = '  +9.59%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.897'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.47%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '308.64'
$c.ClearFormats()
$ws.Range("E38").Value = '  +5.13%  '
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("E40").Value = '  +1.88%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.644'
$c.ClearFormats()
$ws.Range("E41").Value = '  +3.64%  '
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("E43").Value = '  +1.65%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.02%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '19.82'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +1.82%  '
$ws.Range("E47").Value = '  +2.45%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '19.32'
$c.ClearFormats()
$ws.Range("E48").Value = '  +8.41%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("D50").Value = '1.987.01'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +2.41%  '
